# Leandro updated PA datasets
# Update the "Data extracted for ITHIM R" (AE) column to "Yes" and fill in
# the related "Location of files to process original data" (AF) / "Notes"
# (AG) columns for several cities on the "Data" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Mark "Data extracted for ITHIM R" (column AE) as Yes for the rows with
# newly available PA datasets.
$ws.Range("AE5").Value = "Yes"
$ws.Range("AE6").Value = "Yes"
$ws.Range("AE7").Value = "Yes"
$ws.Range("AE8").Value = "Yes"
$ws.Range("AE9").Value = "Yes"
$ws.Range("AE10").Value = "Yes"
$ws.Range("AE11").Value = "Yes"
$ws.Range("AE12").Value = "Yes"

# Fill in "Location of files to process original data" (AF) and "Notes"
# (AG) for the cities with newly processed PA data.

# Row 6 - Colombia / Bogota note
$ws.Range("AG6").Value = "Only LTPA is provided. Sao Paulo occ PA data can be used."

# Row 5 - Brazil / Belo Horizonte
$ws.Range("AF5").Value = "V:\Studies\MOVED\HealthImpact\Data\TIGTHAT\Brazil\Belo Horizonte\Physical activity\Process PA data.R"

# Row 6 - Colombia / Bogota
$ws.Range("AF6").Value = "V:\Studies\MOVED\HealthImpact\Data\TIGTHAT\Colombia\Bogota\Physical activity\Process PA data.R"

# Row 8 - Argentina / Buenos Aires
$ws.Range("AF8").Value = "V:\Studies\MOVED\HealthImpact\Data\TIGTHAT\Argentina\WP3-PA\Process PA data.R"
$ws.Range("AG8").Value = "Only total MVPA can be calculated. However, survey participants indicated in which PA domains they engaged last week (same time of PA data)."

# Row 9 - Chile / Santiago
$ws.Range("AF9").Value = "V:\Studies\MOVED\HealthImpact\Data\TIGTHAT\Chile\Physical activity\Process PA data.R"
$ws.Range("AG9").Value = "Only total MVPA can be calculated. However, survey participants indicated in which PA domains they engaged last week (same time of PA data)."

# Leave the selection on the last cell that was edited.
$ws.Range("AF9").Select()
